$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Distribución de Residuos - first matrix (J1..J4)
Set-TextValue "B5" "0"
Set-TextValue "C6" "0"
Set-TextValue "B7" "0"
Set-TextValue "C8" "0"
Set-TextValue "B9" "0"
Set-TextValue "C9" "0"

# Second matrix (K1..K2) - first block
Set-TextValue "B13" "0.0"
Set-TextValue "B14" "0.0"
Set-TextValue "B17" "0.0"

# Third matrix (K1..K2) - second block
Set-TextValue "B21" "170"
Set-TextValue "B22" "150"
Set-TextValue "C23" "200"
Set-TextValue "C24" "250"
Set-TextValue "B25" "320"
Set-TextValue "C25" "450"

# Distribución de camiones - first matrix (J1..J4)
Set-TextValue "B32" "0.0"
Set-TextValue "C33" "0.0"
Set-TextValue "B34" "0.0"
Set-TextValue "C35" "0.0"
Set-TextValue "B36" "0.0"
Set-TextValue "C36" "0.0"

# Second matrix (K1..K2)
Set-TextValue "B40" "0.0"
Set-TextValue "B41" "0.0"
Set-TextValue "B44" "0.0"

# Third matrix (K1..K2)
Set-TextValue "B48" "195.0"
Set-TextValue "B49" "172.0"
Set-TextValue "C50" "230.0"
Set-TextValue "C51" "287.0"
Set-TextValue "B52" "367.0"
Set-TextValue "C52" "517.0"

# Summary text rows
$ws.Range("A58").Value = "Cantidad de GEI en kg de CO2 Equivalente para el subsistema I-J: 0.0"
$ws.Range("A60").Value = "Cantidad de GEI en kg de CO2 Equivalente para el subsistema J-K: 0.0"
$ws.Range("A62").Value = "Cantidad de GEI en kg de CO2 Equivalente para el subsistema I-K: 48473.0"
$ws.Range("A64").Value = "Cantidad de GEI en kg de CO2 Equivalente para todo el sistema: 48473.0"
$ws.Range("A66").Value = "Costo total de todo el sistema: [8381.]"
$ws.Range("A68").Value = "Valor final de costo de inversión de ET: 4.657541095890411"
$ws.Range("A70").Value = "Valor final de costo operativo de ET: 691.43"
$ws.Range("A72").Value = "Valor final de costo de inversión de CA: 7.934468955701832"
$ws.Range("A74").Value = "Valor final de costo operativo de CA: 11.67591104634009"
